# Auto-generated edit script applying numeric updates per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1160.548
$ws.Range("I15").Value = 1160.548
$ws.Range("K15").Value = 3481.644
$ws.Range("M15").Value = -3312.644
$ws.Range("H108").Value = 28020.25
$ws.Range("J108").Value = 28020.25
$ws.Range("L108").Value = 28020.25
$ws.Range("N108").Value = -35700.25
$ws.Range("H120").Value = 49644
$ws.Range("J120").Value = 49644
$ws.Range("L120").Value = 49644
$ws.Range("N120").Value = -59320
$ws.Range("H126").Value = 46956
$ws.Range("J126").Value = 46956
$ws.Range("L126").Value = 46956
$ws.Range("N126").Value = -56836
$ws.Range("H128").Value = 46672
$ws.Range("J128").Value = 46672
$ws.Range("L128").Value = 46672
$ws.Range("N128").Value = -56632
$ws.Range("H130").Value = 48832
$ws.Range("J130").Value = 48832
$ws.Range("L130").Value = 48832
$ws.Range("N130").Value = -58872

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 53998.4
$ws.Range("J80").Value = 53998.4
$ws.Range("L80").Value = 53998.4
$ws.Range("N80").Value = -55994.4
$ws.Range("H83").Value = 53998.4
$ws.Range("J83").Value = 53998.4
$ws.Range("L83").Value = 161995.2
$ws.Range("N83").Value = -171979.2
$ws.Range("H107").Value = 45220
$ws.Range("J107").Value = 45220
$ws.Range("L107").Value = 45220
$ws.Range("N107").Value = -52900
$ws.Range("H111").Value = 48496
$ws.Range("J111").Value = 48496
$ws.Range("L111").Value = 48496
$ws.Range("N111").Value = -56676
$ws.Range("H119").Value = 52690
$ws.Range("J119").Value = 52690
$ws.Range("L119").Value = 52690
$ws.Range("N119").Value = -62366
$ws.Range("H120").Value = 43296
$ws.Range("J120").Value = 43296
$ws.Range("L120").Value = 43296
$ws.Range("N120").Value = -52972
$ws.Range("H121").Value = 29277.445
$ws.Range("J121").Value = 29277.445
$ws.Range("L121").Value = 29277.445
$ws.Range("N121").Value = -32771.445
$ws.Range("H123").Value = 35610.5
$ws.Range("J123").Value = 35610.5
$ws.Range("L123").Value = 35610.5
$ws.Range("N123").Value = -45410.5
$ws.Range("H125").Value = 50707
$ws.Range("J125").Value = 50707
$ws.Range("L125").Value = 50707
$ws.Range("N125").Value = -60547
$ws.Range("H128").Value = 48025
$ws.Range("J128").Value = 48025
$ws.Range("L128").Value = 48025
$ws.Range("N128").Value = -57985
$ws.Range("H130").Value = 41580.25
$ws.Range("J130").Value = 41580.25
$ws.Range("L130").Value = 41580.25
$ws.Range("N130").Value = -51620.25
$ws.Range("H131").Value = 51695
$ws.Range("J131").Value = 51695
$ws.Range("L131").Value = 51695
$ws.Range("N131").Value = -61775

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 47992
$ws.Range("J124").Value = 47992
$ws.Range("L124").Value = 47992
$ws.Range("N124").Value = -57812
$ws.Range("H125").Value = 50566
$ws.Range("J125").Value = 50566
$ws.Range("L125").Value = 50566
$ws.Range("N125").Value = -60406
$ws.Range("H130").Value = 38494
$ws.Range("J130").Value = 38494
$ws.Range("L130").Value = 38494
$ws.Range("N130").Value = -48534

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 305
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 305
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 305
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -583
$ws.Range("H20").Value = 43481
$ws.Range("J20").Value = 43481
$ws.Range("L20").Value = 43481
$ws.Range("N20").Value = -43953
$ws.Range("H30").Value = 43481
$ws.Range("J30").Value = 43481
$ws.Range("L30").Value = 43481
$ws.Range("N30").Value = -43663
$ws.Range("H100").Value = 25377
$ws.Range("J100").Value = 31754
$ws.Range("L100").Value = 31754
$ws.Range("N100").Value = -33918
$ws.Range("H128").Value = 43481
$ws.Range("J128").Value = 43481
$ws.Range("L128").Value = 43481
$ws.Range("N128").Value = -53441

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2680.4348
$ws.Range("J11").Value = 862.5
$ws.Range("L11").Value = 2587.5
$ws.Range("N11").Value = -2867.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4752.5
$ws.Range("I70").Value = 4824.6875
$ws.Range("J70").Value = 4175
$ws.Range("K70").Value = 4824.6875
$ws.Range("L70").Value = 4175
$ws.Range("M70").Value = -4554.6875
$ws.Range("N70").Value = -4715
$ws.Range("H73").Value = 4752.5
$ws.Range("I73").Value = 4824.6875
$ws.Range("J73").Value = 4175
$ws.Range("K73").Value = 4824.6875
$ws.Range("L73").Value = 4175
$ws.Range("M73").Value = -3888.6875
$ws.Range("N73").Value = -6047
$ws.Range("H104").Value = 32268.834
$ws.Range("J104").Value = 32268.834
$ws.Range("L104").Value = 32268.834
$ws.Range("N104").Value = -39256.834
$ws.Range("H106").Value = 29968
$ws.Range("J106").Value = 29968
$ws.Range("L106").Value = 29968
$ws.Range("N106").Value = -32492
$ws.Range("H110").Value = 31276.8
$ws.Range("J110").Value = 31276.8
$ws.Range("L110").Value = 31276.8
$ws.Range("N110").Value = -39456.8
$ws.Range("H130").Value = 46417.668
$ws.Range("J130").Value = 46417.668
$ws.Range("L130").Value = 46417.668
$ws.Range("N130").Value = -56457.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 3000
$ws.Range("J10").Value = 3000
$ws.Range("L10").Value = 3000
$ws.Range("N10").Value = -3280
$ws.Range("H36").Value = 48707
$ws.Range("J36").Value = 48707
$ws.Range("L36").Value = 48707
$ws.Range("N36").Value = -49831
$ws.Range("H108").Value = 48622
$ws.Range("J108").Value = 48622
$ws.Range("L108").Value = 48622
$ws.Range("N108").Value = -56302
$ws.Range("H111").Value = 43938.5
$ws.Range("J111").Value = 43938.5
$ws.Range("L111").Value = 43938.5
$ws.Range("N111").Value = -52118.5
$ws.Range("H121").Value = 19227.666
$ws.Range("J121").Value = 19227.666
$ws.Range("L121").Value = 19227.666
$ws.Range("N121").Value = -22721.666
$ws.Range("H124").Value = 42140.332
$ws.Range("J124").Value = 42140.332
$ws.Range("L124").Value = 42140.332
$ws.Range("N124").Value = -51960.332
$ws.Range("H127").Value = 42224
$ws.Range("J127").Value = 42224
$ws.Range("L127").Value = 42224
$ws.Range("N127").Value = -52144
$ws.Range("H128").Value = 35210.5
$ws.Range("J128").Value = 35210.5
$ws.Range("L128").Value = 35210.5
$ws.Range("N128").Value = -45170.5
$ws.Range("H130").Value = 38083.332
$ws.Range("J130").Value = 38083.332
$ws.Range("L130").Value = 38083.332
$ws.Range("N130").Value = -48123.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 47210
$ws.Range("J16").Value = 47210
$ws.Range("L16").Value = 47210
$ws.Range("N16").Value = -47794
$ws.Range("H119").Value = 46663.332
$ws.Range("J119").Value = 46663.332
$ws.Range("L119").Value = 46663.332
$ws.Range("N119").Value = -56339.332
$ws.Range("H131").Value = 55515.75
$ws.Range("J131").Value = 55515.75
$ws.Range("L131").Value = 55515.75
$ws.Range("N131").Value = -65595.75
